# Applies the cryptos list price/volume update described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '60.639.29'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +2.95%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.727.58'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +3.86%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.16%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '525.65'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +0.88%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '145.43'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +0.12%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.997'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -0.08%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.576'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +0.86%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '2.726.03'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +3.32%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '6.80'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +7.55%  '
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +1.18%  '
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +1.16%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '3.183.03'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +3.02%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '60.594.29'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '21.26'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +1.91%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.720.74'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +3.33%  '
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +0.97%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '344.63'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -0.13%  '
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +0.41%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '10.61'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +3.92%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.48'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +5.46%  '
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -0.13%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '63.33'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.420'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +1.17%  '
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +1.87%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.995'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -0.04%  '
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +2.47%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.27'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +2.35%  '
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +9.30%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.998'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -0.07%  '
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +1.73%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '19.01'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +0.66%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '149.78'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -0.50%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '4.27'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +7.21%  '
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +7.79%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.941'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -4.01%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.876'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +4.26%  '
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +7.21%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '37.11'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +1.25%  '
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +0.42%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '280.33'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +0.89%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '20.15'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +3.44%  '
$ws.Range('B44').Value = 'FirstDigitalUSD'
$ws.Range('C44').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.998'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +0.27%  '
$ws.Range('B45').Value = 'Mantle'
$ws.Range('C45').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.611'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +0.52%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.143.56'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +7.77%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0985'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +0.05%  '
$ws.Range('B48').Value = 'Hedera'
$ws.Range('C48').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0538'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +2.92%  '
$ws.Range('B49').Value = 'RenderToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '4.83'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +4.36%  '
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +2.36%  '
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +1.59%  '
